$d = $word.ActiveDocument

# Update the date line (first paragraph)
$d.Content.Find.Execute("2025-08-15 Friday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2025-08-16 Saturday", 2)

# Update the division problems table (Table 1), by explicit cell address so the
# two identical "37÷6=" cells are disambiguated by position.
$t = $d.Tables(1)

$values = @{
    1 = @("44÷7=", "24÷8=", "88÷3=", "12÷6=", "56÷4=")
    5 = @("26÷5=", "56÷9=", "30÷3=", "96÷3=", "70÷3=")
    9 = @("74÷6=", "73÷2=", "75÷3=", "12÷8=", "35÷7=")
    13 = @("58÷6=", "73÷5=", "10÷3=", "69÷7=", "90÷9=")
    17 = @("62÷3=", "63÷2=", "61÷6=", "24÷9=", "40÷3=")
}

foreach ($row in $values.Keys) {
    $cols = $values[$row]
    for ($c = 1; $c -le 5; $c++) {
        $cell = $t.Cell($row, $c)
        $rng = $cell.Range
        $rng.MoveEnd(1, -1) | Out-Null
        $rng.Text = $cols[$c - 1]
    }
}
